# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on cells whose new values look numeric,
# so Excel stores them as text (matching original inlineStr cells)
# instead of auto-converting them to numbers. Applied per-cell (not as
# a combined multi-area range) since multi-area NumberFormat only
# reliably affects the first area in this runtime.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "29.283.44"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.872.27"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D5").Value = "0.7094"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "241.53"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.07799"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "0.3094"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "25.01"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "0.08397"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.871.23"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "5.231"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "0.7106"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "91.12"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "29.301.20"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "0.000008165"
$ws.Range("E18").Value = "  +4.22%  "
$ws.Range("D19").Value = "239.93"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "13.21"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "2.126.48"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D23").Value = "7.742"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "0.1598"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").Value = "162.77"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "8.993"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "18.46"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "4.386"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "1.297"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "4.291"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "0.05374"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "1.944"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "0.7479"
$ws.Range("E36").Value = "  -5.80%  "
$ws.Range("D37").Value = "2.695"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "1.236.26"
$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "6.508"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").Value = "0.8895"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "72.28"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "108.35"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "2.020.18"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "0.5196"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000123"
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.788"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "9.403"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "0.4308"
$ws.Range("E51").Value = "  +0.32%  "

# Restore default (no explicit) formatting on those cells
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
